$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (number of people interested) column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 1887
$wsExpo.Range("F4").Value = 823
$wsExpo.Range("F5").Value = 738
$wsExpo.Range("F6").Value = 236

# Sheet "全部类型" (All types) - update same "想去人数" column F (different row mapping)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 1887
$wsAll.Range("F5").Value = 823
$wsAll.Range("F6").Value = 738
$wsAll.Range("F7").Value = 236
